# Weekly Fruta/Hortalizas update: insert one new price record row
# before the existing row 1000 (shifting all subsequent rows down by one),
# and populate it with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 1000, pushing former rows 1000..1086
# down to 1001..1087 (this also extends the used range / dimension
# automatically to A1:T1087).
$ws.Rows.Item(1000).Insert()

# Populate the newly inserted row with this week's record.
$ws.Cells.Item(1000, 1).Value  = 3
$ws.Cells.Item(1000, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(1000, 3).Value  = "Coquimbo"
$ws.Cells.Item(1000, 4).Value  = 45132
$ws.Cells.Item(1000, 5).Value  = 5
$ws.Cells.Item(1000, 6).Value  = "Fruta"
$ws.Cells.Item(1000, 7).Value  = 100108
$ws.Cells.Item(1000, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(1000, 9).Value  = 100108005
$ws.Cells.Item(1000, 10).Value = "Piña"
$ws.Cells.Item(1000, 11).Value = "Caramelo"
$ws.Cells.Item(1000, 12).Value = "Primera"
$ws.Cells.Item(1000, 13).Value = 108
$ws.Cells.Item(1000, 14).Value = 25000
$ws.Cells.Item(1000, 15).Value = 25000
$ws.Cells.Item(1000, 16).Value = 25000
$ws.Cells.Item(1000, 17).Value = "`$/caja 12 unidades"
$ws.Cells.Item(1000, 18).Value = "Ecuador"
$ws.Cells.Item(1000, 19).Value = 2083
$ws.Cells.Item(1000, 20).Value = 12
